$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - update the sample-size header values for columns B:E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 - B2, D2, E2 were deleted (cleared); C2 got a new value
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 6.0407726312651349
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3 - values tweaked
$ws.Range("B3").Value = 6.3751365426387139
$ws.Range("C3").Value = 6.6796629283111173
$ws.Range("D3").Value = 8.099961900979336
$ws.Range("E3").Value = 3.3055796374253683

# Selection now only covers B1:E3 instead of the full B1:AY3
$ws.Range("B1:E3").Select()
